# Extend the "Studenti" sheet with a new "Courses" column (H), populate the
# header + first student's courses, size the new column, and scroll/select
# so column H is visible with H2 as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (H1) and sample data (H2)
$ws.Range("H1").Value = "Courses"
$ws.Range("H2").Value = ".NET, Android, Cloud Computing"

# Widen the new column to fit the course list text
$ws.Range("H1").ColumnWidth = 50.6666666666667

# Scroll the view so column B becomes the left-most visible column, then
# select H2 (matches the saved view state / active cell)
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H2").Select()
